$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the libraryProtocol value in column K (E7760 -> E7420) for all data rows
$kRange = $ws.Range("K2:K27")
$kRange.Value = "E7420"

# 2. Apply a new font look to column K (size 11, explicit black color), no wrap
$kRange.Font.Size = 11
$kRange.Font.Color = 0
$kRange.WrapText = $false

# 3. Replace the literal boolean in column L with a live formula
for ($row = 2; $row -le 27; $row++) {
    $ws.Cells.Item($row, 12).Formula = "=FALSE()"
}

# 4. Leave the final selection on column K (mirrors the last-edited range)
$kRange.Select() | Out-Null
